$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Sheet 1: "Rushing" - update cumulative attempt stats after the
# Wild Card round game (no new rows, only value changes)
# -----------------------------------------------------------------
$rushing = $wb.Worksheets.Item(1)

$rushing.Range("C2").Value = 11   # R.Wilson 1DATT
$rushing.Range("E2").Value = 17   # R.Wilson 3DATT
$rushing.Range("F2").Value = 5    # R.Wilson RZATT

$rushing.Range("C3").Value = 67   # R.Penny 1DATT
$rushing.Range("D3").Value = 41   # R.Penny 2DATT
$rushing.Range("F3").Value = 17   # R.Penny RZATT

$rushing.Range("C5").Value = 4    # T.Homer 1DATT

$rushing.Range("C6").Value = 5    # D.Dallas 1DATT

$rushing.Range("C9").Value = 3    # T.Lockett 1DATT

# -----------------------------------------------------------------
# Sheet 2: "Receiving" - log the Wild Card round: T.Homer had his
# first receiving action this round, so a new row is inserted for
# him (row 4), shifting everyone below down by one. Several other
# players' cumulative totals also grow, and C.Parkinson (now on the
# last row) gets his post-game totals too.
# -----------------------------------------------------------------
$receiving = $wb.Worksheets.Item(2)

# Copy the row-13 formatting down into the brand-new row 14 so the
# new last row matches the look (bold/bordered index column) of the
# rest of the table without touching the shared style table.
$receiving.Range("A13").Copy()
$receiving.Range("A14").PasteSpecial(-4122)  # xlPasteFormats

# Shift rows 5-13 down to 6-14, writing each row's final content
# (this naturally makes room for the inserted T.Homer row at 4).
$receiving.Cells.Item(14, 1).Value = 12
$receiving.Cells.Item(14, 2).Value = "C.Parkinson"
$receiving.Cells.Item(14, 3).Value = 6
$receiving.Cells.Item(14, 4).Value = 5
$receiving.Cells.Item(14, 5).Value = 2
$receiving.Cells.Item(14, 6).Value = 0
$receiving.Cells.Item(14, 7).Value = 3
$receiving.Cells.Item(14, 8).Value = 1

$receiving.Cells.Item(13, 1).Value = 11
$receiving.Cells.Item(13, 2).Value = "W.Dissly"
$receiving.Cells.Item(13, 3).Value = 21
$receiving.Cells.Item(13, 4).Value = 16
$receiving.Cells.Item(13, 5).Value = 5
$receiving.Cells.Item(13, 6).Value = 5
$receiving.Cells.Item(13, 7).Value = 3
$receiving.Cells.Item(13, 8).Value = 2

$receiving.Cells.Item(12, 1).Value = 10
$receiving.Cells.Item(12, 2).Value = "G.Everett"
$receiving.Cells.Item(12, 3).Value = 56
$receiving.Cells.Item(12, 4).Value = 43
$receiving.Cells.Item(12, 5).Value = 7
$receiving.Cells.Item(12, 6).Value = 5
$receiving.Cells.Item(12, 7).Value = 8
$receiving.Cells.Item(12, 8).Value = 4

$receiving.Cells.Item(11, 1).Value = 9
$receiving.Cells.Item(11, 2).Value = "P.Hart"
$receiving.Cells.Item(11, 3).Value = 10
$receiving.Cells.Item(11, 4).Value = 7
$receiving.Cells.Item(11, 5).Value = 2
$receiving.Cells.Item(11, 6).Value = 0
$receiving.Cells.Item(11, 7).Value = 1
$receiving.Cells.Item(11, 8).Value = 0

$receiving.Cells.Item(10, 1).Value = 8
$receiving.Cells.Item(10, 2).Value = "D.Eskridge"
$receiving.Cells.Item(10, 3).Value = 13
$receiving.Cells.Item(10, 4).Value = 10
$receiving.Cells.Item(10, 5).Value = 7
$receiving.Cells.Item(10, 6).Value = 0
$receiving.Cells.Item(10, 7).Value = 4
$receiving.Cells.Item(10, 8).Value = 2

$receiving.Cells.Item(9, 1).Value = 7
$receiving.Cells.Item(9, 2).Value = "F.Swain"
$receiving.Cells.Item(9, 3).Value = 31
$receiving.Cells.Item(9, 4).Value = 19
$receiving.Cells.Item(9, 5).Value = 9
$receiving.Cells.Item(9, 6).Value = 6
$receiving.Cells.Item(9, 7).Value = 4
$receiving.Cells.Item(9, 8).Value = 1

$receiving.Cells.Item(8, 1).Value = 6
$receiving.Cells.Item(8, 2).Value = "T.Lockett"
$receiving.Cells.Item(8, 3).Value = 70
$receiving.Cells.Item(8, 4).Value = 54
$receiving.Cells.Item(8, 5).Value = 46
$receiving.Cells.Item(8, 6).Value = 27
$receiving.Cells.Item(8, 7).Value = 9
$receiving.Cells.Item(8, 8).Value = 5

$receiving.Cells.Item(7, 1).Value = 5
$receiving.Cells.Item(7, 2).Value = "D.Metcalf"
$receiving.Cells.Item(7, 3).Value = 91
$receiving.Cells.Item(7, 4).Value = 64
$receiving.Cells.Item(7, 5).Value = 38
$receiving.Cells.Item(7, 6).Value = 11
$receiving.Cells.Item(7, 7).Value = 20
$receiving.Cells.Item(7, 8).Value = 13

$receiving.Cells.Item(6, 1).Value = 4
$receiving.Cells.Item(6, 2).Value = "N.Bellore"
$receiving.Cells.Item(6, 3).Value = 1
$receiving.Cells.Item(6, 4).Value = 1
$receiving.Cells.Item(6, 5).Value = 0
$receiving.Cells.Item(6, 6).Value = 0
$receiving.Cells.Item(6, 7).Value = 0
$receiving.Cells.Item(6, 8).Value = 0

$receiving.Cells.Item(5, 1).Value = 3
$receiving.Cells.Item(5, 2).Value = "D.Dallas"
$receiving.Cells.Item(5, 3).Value = 13
$receiving.Cells.Item(5, 4).Value = 11
$receiving.Cells.Item(5, 5).Value = 0
$receiving.Cells.Item(5, 6).Value = 0
$receiving.Cells.Item(5, 7).Value = 4
$receiving.Cells.Item(5, 8).Value = 3

# Newly logged row for T.Homer (row 4), between A.Collins and D.Dallas
$receiving.Cells.Item(4, 1).Value = 2
$receiving.Cells.Item(4, 2).Value = "T.Homer"
$receiving.Cells.Item(4, 3).Value = 2
$receiving.Cells.Item(4, 4).Value = 0
$receiving.Cells.Item(4, 5).Value = 1
$receiving.Cells.Item(4, 6).Value = 1
$receiving.Cells.Item(4, 7).Value = 0
$receiving.Cells.Item(4, 8).Value = 0
